# Weekly update: insert a new week's Broccoli price report (2 rows) at the
# top of the data block (rows 423-424), pushing the existing history down.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the current first data row of this block (423),
# shifting rows 423:440 down to 425:442. Use the existing rows as the
# source of the insert so number formats (e.g. the date style on column D)
# carry over to the freshly inserted rows.
$ws.Range("A423:R424").EntireRow.Insert()

# New week (2021-11-09, serial 44509) - "Primera" quality row.
$ws.Cells.Item(423, 1).Value = 8
$ws.Cells.Item(423, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(423, 3).Value = "Coquimbo"
$ws.Cells.Item(423, 4).Value = 44509
$ws.Cells.Item(423, 5).Value = 4
$ws.Cells.Item(423, 6).Value = 100112023
$ws.Cells.Item(423, 7).Value = "Brócoli"
$ws.Cells.Item(423, 8).Value = "Sin especificar"
$ws.Cells.Item(423, 9).Value = "Primera"
$ws.Cells.Item(423, 10).Value = 2200
$ws.Cells.Item(423, 11).Value = 600
$ws.Cells.Item(423, 12).Value = 700
$ws.Cells.Item(423, 13).Value = 650
$ws.Cells.Item(423, 14).Value = "`$/unidad"
$ws.Cells.Item(423, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(423, 16).Value = 650
$ws.Cells.Item(423, 17).Value = 1
$ws.Cells.Item(423, 18).Value = "Hortaliza"

# New week (2021-11-09, serial 44509) - "Segunda" quality row.
$ws.Cells.Item(424, 1).Value = 8
$ws.Cells.Item(424, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(424, 3).Value = "Coquimbo"
$ws.Cells.Item(424, 4).Value = 44509
$ws.Cells.Item(424, 5).Value = 4
$ws.Cells.Item(424, 6).Value = 100112023
$ws.Cells.Item(424, 7).Value = "Brócoli"
$ws.Cells.Item(424, 8).Value = "Sin especificar"
$ws.Cells.Item(424, 9).Value = "Segunda"
$ws.Cells.Item(424, 10).Value = 1300
$ws.Cells.Item(424, 11).Value = 500
$ws.Cells.Item(424, 12).Value = 550
$ws.Cells.Item(424, 13).Value = 525
$ws.Cells.Item(424, 14).Value = "`$/unidad"
$ws.Cells.Item(424, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(424, 16).Value = 525
$ws.Cells.Item(424, 17).Value = 1
$ws.Cells.Item(424, 18).Value = "Hortaliza"
